# Corrected excel sheets for application fix issues
#
# On the "Edit Repayment Schedule" sheet, insert a new automation step
# ("waittopageload1" / 2000) right after "validaterepaymentschedule"
# (old row 6) and before "clickonsubmit" (old row 6, new row 7), pushing
# every row below it down by one.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("NewLoanInput")
$wsEdit = $wb.Worksheets.Item("Edit Repayment Schedule")

# Insert a new row 6 (shifts old rows 6-12 down to 7-13).
$wsEdit.Rows.Item(6).Insert()

# Populate the new row: A6 = "waittopageload1", B6 = 2000.
$wsEdit.Range("A6").Value = "waittopageload1"
$wsEdit.Range("B6").Value = 2000

# Match the formatting used by the existing "waittopageload" row (B3),
# which carries the numeric wait-time style.
$wsEdit.Range("B3").Copy()
$wsEdit.Range("B6").PasteSpecial(-4122)

# Reflect the new selection on this sheet (A6:B6), then restore the
# workbook's originally active sheet/cell (NewLoanInput!B2) so the only
# sheet-view change that sticks is the one on "Edit Repayment Schedule".
[void]$wsEdit.Range("A6:B6").Select()
$ws1.Activate()
[void]$ws1.Range("B2").Select()
